# Auto-upload VRF Excel file
# Adds a new "cat" worksheet at the end of the workbook with the standard
# VRF model header row (Outdoor/Indoor Model, Quantity, Serial(s)).

$wb = $excel.ActiveWorkbook

# Remember the sheet that is currently last -- the new sheet will be
# inserted right after it (i.e. appended at the end of the tab strip),
# and it also doubles as the formatting template for the header row that
# every other sheet in this workbook already uses.
$sheetCount = $wb.Worksheets.Count
$templateSheet = $wb.Worksheets.Item($sheetCount)

$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $templateSheet)
$newSheet.Name = "cat"

# Standard header row used by every sheet in this workbook.
$headers = @("Outdoor Model", "Outdoor Quantity", "Outdoor Serial(s)", "Indoor Model", "Indoor Quantity", "Indoor Serial(s)")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $newSheet.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Match the bold / thin-bordered / centered-top header formatting used
# throughout the rest of the workbook by copying it from the template
# sheet's own header row.
$templateSheet.Range("A1:F1").Copy()
$newSheet.Range("A1:F1").PasteSpecial(-4122)

$newSheet.Range("A1").Select()
